# Populates the Tracking Log audit rows (142-256) on "Hoja1" with the
# real audit events that were previously placeholder/blank rows.
# Source data: row -> (Usuario shared-text id, Fecha/Hora serial, Descripcion shared-text id)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Lookup table for the shared text used in the "Usuario" (B) and
# "Descripcion" (E) columns, keyed by the id used in $rowsData below.
$textById = @{
    7 = "riesgos"
    8 = "riesgos ingreso al sistema"
    11 = "riesgos salio del sistema"
    16 = "riesgos dio de alta a Miguel Sandoval"
    17 = "riesgos modifico permisos del usuario riesgos"
    18 = "riesgos modifico permisos del usuario Miguel Sandoval"
    19 = "riesgos dio de alta a coca"
    20 = "riesgos modifico permisos del usuario coca"
    21 = "coca"
    22 = "coca ingreso al sistema"
    23 = "coca salio del sistema"
    24 = "riesgos dio de baja al usuario coca"
    25 = "riesgos habilito al usuario coca"
    26 = "Genero reporte Reporte RC01"
    27 = "Genero reporte Reporte ICAP"
    28 = "Genero reporte Tracking Log"
}

# Compact per-row records: Row|UsuarioId|FechaHoraSerial|DescripcionId
$rowsData = @"
142|7|41802.97175925926|8
143|7|41802.973078703704|8
144|7|41802.98372685185|8
145|7|41802.9852662037|8
146|7|41802.98813657407|8
147|7|41802.99209490741|8
148|7|41802.996886574074|8
149|7|41804.51609953704|8
150|7|41804.524502314816|16
151|7|41804.52474537037|17
152|7|41804.52670138889|18
153|7|41804.526967592596|11
154|7|41804.52890046296|8
155|7|41806.91587962963|8
156|7|41806.91914351852|8
157|7|41806.92680555556|8
158|7|41806.92697916667|11
159|7|41806.92707175926|8
160|7|41806.93303240741|8
161|7|41806.93318287037|11
162|7|41806.93561342593|8
163|7|41806.93587962963|11
164|7|41806.93671296296|8
165|7|41806.94064814815|8
166|7|41806.94326388889|8
167|7|41806.9534375|8
168|7|41806.95380787037|11
169|7|41806.95386574074|8
170|7|41806.95394675926|11
171|7|41806.95758101852|8
172|7|41806.95767361111|11
173|7|41806.95930555555|8
174|7|41806.959444444445|11
175|7|41806.96302083333|8
176|7|41806.964907407404|8
177|7|41806.9653125|19
178|7|41806.96542824074|20
179|7|41806.96545138889|11
180|21|41806.96549768518|22
181|21|41806.96561342593|23
182|7|41806.9656712963|8
183|7|41806.965891203705|24
184|7|41806.96601851852|11
185|7|41806.966261574074|8
186|7|41806.96681712963|11
187|7|41806.96957175926|8
188|7|41806.9696875|24
189|7|41806.9699537037|11
190|7|41806.97740740741|8
191|7|41806.97760416667|24
192|7|41806.977638888886|11
193|7|41806.97872685185|8
194|7|41806.97886574074|24
195|7|41806.97893518519|11
196|21|41806.978993055556|22
197|21|41806.97902777778|23
198|7|41806.97907407407|8
199|7|41806.97917824074|24
200|7|41806.97918981482|11
201|7|41806.98150462963|8
202|7|41806.981828703705|25
203|7|41806.981886574074|25
204|7|41806.982152777775|11
205|7|41806.98974537037|8
206|7|41806.9999537037|8
207|7|41807.005960648145|8
208|7|41807.00616898148|20
209|7|41807.00625|20
210|7|41807.00628472222|11
211|21|41807.00633101852|22
212|21|41807.00646990741|23
213|7|41807.00885416667|8
214|7|41807.00913194445|11
215|7|41807.014918981484|8
216|7|41807.02423611111|8
217|21|41807.027083333334|22
218|7|41807.03873842592|8
219|7|41807.04425925926|8
220|7|41807.04614583333|8
221|21|41807.0540625|22
222|21|41807.05454861111|26
223|21|41807.0553125|27
224|21|41807.056550925925|27
225|21|41807.057974537034|23
226|21|41807.0580787037|22
227|21|41807.05868055556|26
228|7|41807.06298611111|8
229|7|41807.06392361111|26
230|7|41807.06418981482|11
231|7|41807.32502314815|8
232|7|41807.879791666666|8
233|7|41807.88081018518|26
234|21|41807.88452546296|22
235|7|41807.88659722222|8
236|7|41807.88744212963|11
237|7|41807.88753472222|8
238|21|41807.888761574075|22
239|21|41807.907847222225|23
240|7|41807.90828703704|8
241|7|41807.90846064815|11
242|7|41807.908530092594|8
243|7|41807.92271990741|8
244|7|41807.92291666667|11
245|7|41807.928761574076|8
246|7|41807.928923611114|28
247|7|41807.93728009259|8
248|7|41807.93746527778|28
249|7|41807.94180555556|8
250|7|41807.94196759259|28
251|7|41807.94372685185|8
252|7|41807.94385416667|28
253|7|41807.94835648148|8
254|7|41807.94850694444|28
255|7|41807.95048611111|8
256|7|41807.95065972222|28
"@

$lines = $rowsData -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split "\|"
    $rowNum = [int]$parts[0]
    $userId = [int]$parts[1]
    $serial = [double]$parts[2]
    $descId = [int]$parts[3]

    $userText = $textById[$userId]
    $descText = $textById[$descId]

    $ws.Cells.Item($rowNum, 2).Value2 = $userText
    $ws.Cells.Item($rowNum, 3).Value2 = $serial
    $ws.Cells.Item($rowNum, 4).Value2 = $serial
    $ws.Cells.Item($rowNum, 5).Value2 = $descText
}

